$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.004.82"
$ws.Range("E2").Value = "  +0.96%  "
$ws.Range("D3").Value = "2.236.63"
$ws.Range("E3").Value = "  +1.66%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "98.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +17.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "270.90"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.77%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.26%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.643"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +8.17%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "48.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.88%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0949"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +19.31%  "
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.64%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.822"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.63%  "
$ws.Range("D16").Value = "2.249.44"
$ws.Range("E16").Value = "  +2.29%  "
$ws.Range("D17").Value = "2.163.04"
$ws.Range("E17").Value = "  -14.30%  "
$ws.Range("D18").Value = "43.988.77"
$ws.Range("E18").Value = "  +1.10%  "
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.30"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.38%  "
$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.11%  "
$ws.Range("E27").Value = "  +12.97%  "
$ws.Range("B28").Value = "WEMIXToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.97%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0924"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.28%  "
$ws.Range("B35").Value = "Stellar"
$ws.Range("C35").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.125"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.114"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0352"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.36"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.25%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +29.99%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.251"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +25.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "63.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.43"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.81%  "
$ws.Range("E45").Value = "  +4.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.25%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.16"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.82%  "
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.439"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.05%  "
